$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($row, $values) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}

# Row 2: Target cluster ECs, refreshed TPM-derived values
Set-RowValues 2 @(
    "ECs", "Nxph3", "Nrxn2", "ECs", 2, 0.6666666666666666, 0.4198326666666667, 1.259498,
    1, 1, 2, 0.6666666666666666, 0.136464, 0.409392, 0.2657380151046518,
    0.2657380151046518, 0.057292045024, 0.515628405216, 0.2657380151046518, 0.2657380151046518
)

# Row 3: Target cluster FAPs, refreshed TPM-derived values
Set-RowValues 3 @(
    "ECs", "Nxph3", "Nrxn2", "FAPs", 2, 0.6666666666666666, 0.4198326666666667, 1.259498,
    1, 1, 3, 1, 0.3195506666666667, 0.9586520000000001, 0.6222649188457632,
    0.6222649188457632, 0.1341578085217778, 1.207420276696, 0.6222649188457632, 0.6222649188457632
)

# Row 4: Target cluster MuSCs (moved from row 2), refreshed TPM-derived values
Set-RowValues 4 @(
    "ECs", "Nxph3", "Nrxn2", "MuSCs", 2, 0.6666666666666666, 0.4198326666666667, 1.259498,
    1, 1, 1, 0.3333333333333333, 0.01800166666666667, 0.054005, 0.03505486552186345,
    0.03505486552186345, 0.007557687721111111, 0.06801918948999999, 0.03505486552186345, 0.03505486552186345
)

# Row 5 (new): Target cluster Neutrophils
Set-RowValues 5 @(
    "ECs", "Nxph3", "Nrxn2", "Neutrophils", 2, 0.6666666666666666, 0.4198326666666667, 1.259498,
    1, 1, 1, 0.3333333333333333, 0.039512, 0.118536, 0.0769422005277216,
    0.0769422005277216, 0.01658842832533333, 0.149295854928, 0.0769422005277216, 0.0769422005277216
)

$wb.Save()
